# Apply updated cryptocurrency Price (D) and Volume(1h) (E) values.
# Values are written via the COM Range.Value setter. Where the new
# Price text looks like a plain number, a leading apostrophe is used
# (exactly as typing `'0.9992` into Excel would) so the cell keeps
# storing literal text instead of being auto-converted to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.738.07'
$ws.Range("E2").Value = '  +6.85%  '
$ws.Range("D3").Value = '1.812.36'
$ws.Range("E3").Value = '  +4.94%  '
$ws.Range("D4").Value = '''0.9992'
$ws.Range("D5").Value = '''251.01'
$ws.Range("E5").Value = '  +3.73%  '
$ws.Range("D6").Value = '''0.9992'
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("D7").Value = '''0.4974'
$ws.Range("E7").Value = '  +1.58%  '
$ws.Range("D8").Value = '''0.2786'
$ws.Range("E8").Value = '  +7.59%  '
$ws.Range("D9").Value = '''0.06383'
$ws.Range("E9").Value = '  +2.71%  '
$ws.Range("D10").Value = '1.809.93'
$ws.Range("E10").Value = '  +4.63%  '
$ws.Range("D11").Value = '''16.73'
$ws.Range("E11").Value = '  +4.62%  '
$ws.Range("D12").Value = '''0.07108'
$ws.Range("E12").Value = '  +2.97%  '
$ws.Range("D13").Value = '''0.6479'
$ws.Range("E13").Value = '  +6.38%  '
$ws.Range("D14").Value = '''4.703'
$ws.Range("E14").Value = '  +4.86%  '
$ws.Range("D15").Value = '''81.78'
$ws.Range("E15").Value = '  +5.87%  '
$ws.Range("D16").Value = '28.707.16'
$ws.Range("E16").Value = '  +6.82%  '
$ws.Range("D17").Value = '''0.9991'
$ws.Range("E17").Value = '  +0.09%  '
$ws.Range("D18").Value = '''0.000007382'
$ws.Range("E18").Value = '  +2.78%  '
$ws.Range("D19").Value = '''0.9988'
$ws.Range("E19").Value = '  +0.14%  '
$ws.Range("D20").Value = '''12.28'
$ws.Range("E20").Value = '  +7.30%  '
$ws.Range("D21").Value = '2.037.88'
$ws.Range("D22").Value = '''4.616'
$ws.Range("E22").Value = '  +4.28%  '
$ws.Range("D23").Value = '''8.887'
$ws.Range("E23").Value = '  +3.65%  '
$ws.Range("D24").Value = '''5.327'
$ws.Range("E24").Value = '  +4.44%  '
$ws.Range("D25").Value = '''143.06'
$ws.Range("E25").Value = '  +3.19%  '
$ws.Range("D26").Value = '''16.02'
$ws.Range("E26").Value = '  +4.66%  '
$ws.Range("D27").Value = '''1.878'
$ws.Range("E27").Value = '  +4.67%  '
$ws.Range("D28").Value = '''112.67'
$ws.Range("E28").Value = '  +6.22%  '
$ws.Range("D29").Value = '''1.391'
$ws.Range("E29").Value = '  +0.78%  '
$ws.Range("D30").Value = '''4.175'
$ws.Range("E30").Value = '  +5.77%  '
$ws.Range("D31").Value = '''0.08371'
$ws.Range("E31").Value = '  +4.72%  '
$ws.Range("D32").Value = '''3.840'
$ws.Range("E32").Value = '  +4.12%  '
$ws.Range("D33").Value = '''0.04971'
$ws.Range("E33").Value = '  +9.69%  '
$ws.Range("E34").Value = '  +7.89%  '
$ws.Range("D35").Value = '''0.6766'
$ws.Range("E35").Value = '  +8.14%  '
$ws.Range("D36").Value = '''2.662'
$ws.Range("E36").Value = '  +2.51%  '
$ws.Range("D37").Value = '''2.695'
$ws.Range("E37").Value = '  +9.71%  '
$ws.Range("D38").Value = '''0.9594'
$ws.Range("E38").Value = '  +2.29%  '
$ws.Range("D39").Value = '''2.140'
$ws.Range("E39").Value = '  +4.13%  '
$ws.Range("D40").Value = '''0.01592'
$ws.Range("E40").Value = '  +5.97%  '
$ws.Range("D41").Value = '''5.976'
$ws.Range("E41").Value = '  +5.72%  '
$ws.Range("D42").Value = '''0.9999'
$ws.Range("E42").Value = '  +0.19%  '
$ws.Range("D43").Value = '''101.19'
$ws.Range("E43").Value = '  +1.77%  '
$ws.Range("E44").Value = '  +6.72%  '
$ws.Range("D45").Value = '''7.191'
$ws.Range("E45").Value = '  +4.41%  '
$ws.Range("E46").Value = '  +5.25%  '
$ws.Range("D47").Value = '''0.05499'
$ws.Range("E47").Value = '  +1.95%  '
$ws.Range("D48").Value = '''8.173'
$ws.Range("E48").Value = '  +3.37%  '
$ws.Range("D49").Value = '''31.44'
$ws.Range("E49").Value = '  +4.19%  '
$ws.Range("D50").Value = '''0.3629'
$ws.Range("E50").Value = '  +7.38%  '
$ws.Range("D51").Value = '''1.302'
$ws.Range("E51").Value = '  +5.33%  '
